$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = New-Object 'object[,]' 24,1
$colB[0,0] = 1.02
$colB[1,0] = 1.02
$colB[2,0] = 1.02
$colB[3,0] = 1.02
$colB[4,0] = 1.02
$colB[5,0] = 1.02
$colB[6,0] = 1.02
$colB[7,0] = 1.02
$colB[8,0] = 1.02
$colB[9,0] = 1.02
$colB[10,0] = 1.02
$colB[11,0] = 1.02
$colB[12,0] = 1.02
$colB[13,0] = 1.02
$colB[14,0] = 1.02
$colB[15,0] = 1.02
$colB[16,0] = 1.02
$colB[17,0] = 1.02
$colB[18,0] = 1.02
$colB[19,0] = 1.02
$colB[20,0] = 1.02
$colB[21,0] = 1.02
$colB[22,0] = 1.02
$colB[23,0] = 1.02
$ws.Range("B2:B25").Value = $colB

$colC = New-Object 'object[,]' 24,1
$colC[0,0] = 1.037174652730592
$colC[1,0] = 1.038054295861997
$colC[2,0] = 1.038623885310595
$colC[3,0] = 1.038863435537514
$colC[4,0] = 1.038903662597334
$colC[5,0] = 1.038627085821487
$colC[6,0] = 1.037471848323503
$colC[7,0] = 1.035439304997625
$colC[8,0] = 1.034086459169527
$colC[9,0] = 1.033501195698953
$colC[10,0] = 1.03328388331554
$colC[11,0] = 1.033330493895324
$colC[12,0] = 1.033483230939566
$colC[13,0] = 1.033577347981654
$colC[14,0] = 1.034125312402989
$colC[15,0] = 1.034469178167397
$colC[16,0] = 1.034669800081639
$colC[17,0] = 1.034738215553882
$colC[18,0] = 1.034432279358409
$colC[19,0] = 1.033438251453894
$colC[20,0] = 1.032813733190112
$colC[21,0] = 1.033144757498677
$colC[22,0] = 1.034448952184954
$colC[23,0] = 1.035964387020788
$ws.Range("C2:C25").Value = $colC

$colD = New-Object 'object[,]' 24,1
$colD[0,0] = 1.039345150821084
$colD[1,0] = 1.039992116374041
$colD[2,0] = 1.040410971569566
$colD[3,0] = 1.04058711067115
$colD[4,0] = 1.040616688227943
$colD[5,0] = 1.040413324944993
$colD[6,0] = 1.039563748299359
$colD[7,0] = 1.038068478633945
$colD[8,0] = 1.037072929043329
$colD[9,0] = 1.036642171859087
$colD[10,0] = 1.036482219428598
$colD[11,0] = 1.036516527482143
$colD[12,0] = 1.036628949113091
$colD[13,0] = 1.036698222457359
$colD[14,0] = 1.037101523901912
$colD[15,0] = 1.037354591632898
$colD[16,0] = 1.037502232758793
$colD[17,0] = 1.037552579778686
$colD[18,0] = 1.037327436642795
$colD[19,0] = 1.036595842353146
$colD[20,0] = 1.036136149840036
$colD[21,0] = 1.036379813551689
$colD[22,0] = 1.037339706718344
$colD[23,0] = 1.038454819786583
$ws.Range("D2:D25").Value = $colD

$colE = New-Object 'object[,]' 24,1
$colE[0,0] = 1.044846267074777
$colE[1,0] = 1.045662683957636
$colE[2,0] = 1.04619175919699
$colE[3,0] = 1.046414372159535
$colE[4,0] = 1.046451760936924
$colE[5,0] = 1.046194733018695
$colE[6,0] = 1.04512201263919
$colE[7,0] = 1.043237928755869
$colE[8,0] = 1.041986123100413
$colE[9,0] = 1.041445104562645
$colE[10,0] = 1.041244301190728
$colE[11,0] = 1.041287367157046
$colE[12,0] = 1.041428502922486
$colE[13,0] = 1.041515481910428
$colE[14,0] = 1.042022050372244
$colE[15,0] = 1.042340081963367
$colE[16,0] = 1.042525682874575
$colE[17,0] = 1.042588984643961
$colE[18,0] = 1.042305949987785
$colE[19,0] = 1.041386937666245
$colE[20,0] = 1.040810015690758
$colE[21,0] = 1.04111576736033
$colE[22,0] = 1.042321372454288
$colE[23,0] = 1.043724267111976
$ws.Range("E2:E25").Value = $colE

$colF = New-Object 'object[,]' 24,1
$colF[0,0] = 1.052998016541884
$colF[1,0] = 1.053975934325667
$colF[2,0] = 1.054609856701186
$colF[3,0] = 1.054876630050586
$colF[4,0] = 1.054921438397864
$colF[5,0] = 1.05461342027138
$colF[6,0] = 1.053328270601946
$colF[7,0] = 1.05107250361547
$colF[8,0] = 1.049574687293531
$colF[9,0] = 1.048927565109513
$colF[10,0] = 1.04868741329845
$colF[11,0] = 1.048738916779681
$colF[12,0] = 1.048907709626312
$colF[13,0] = 1.049011737426377
$colF[14,0] = 1.049617665175447
$colF[15,0] = 1.049998134819793
$colF[16,0] = 1.050220195373043
$colF[17,0] = 1.050295935821921
$colF[18,0] = 1.049957299634496
$colF[19,0] = 1.048857998305008
$colF[20,0] = 1.048168087007256
$colF[21,0] = 1.048533701808992
$colF[22,0] = 1.049975750867528
$colF[23,0] = 1.051654617580936
$ws.Range("F2:F25").Value = $colF

$colI = New-Object 'object[,]' 24,1
$colI[0,0] = 1.039510812905085
$colI[1,0] = 1.039731913920884
$colI[2,0] = 1.039873779529956
$colI[3,0] = 1.039933131484264
$colI[4,0] = 1.039943080009421
$colI[5,0] = 1.039874573727554
$colI[6,0] = 1.039585783393737
$colI[7,0] = 1.039067726477224
$colI[8,0] = 1.038716230146761
$colI[9,0] = 1.038562584177112
$colI[10,0] = 1.038505296681528
$colI[11,0] = 1.038517594832068
$colI[12,0] = 1.038557853188385
$colI[13,0] = 1.038582629016477
$colI[14,0] = 1.038726396720382
$colI[15,0] = 1.038816191806088
$colI[16,0] = 1.038868428222873
$colI[17,0] = 1.038886215802909
$colI[18,0] = 1.038806572071416
$colI[19,0] = 1.03854600407812
$colI[20,0] = 1.038380922127256
$colI[21,0] = 1.038468553694209
$colI[22,0] = 1.03881091924665
$colI[23,0] = 1.039202739219382
$ws.Range("I2:I25").Value = $colI

$colJ = New-Object 'object[,]' 24,1
$colJ[0,0] = 1.042278957846756
$colJ[1,0] = 1.042803222087647
$colJ[2,0] = 1.043142209178984
$colJ[3,0] = 1.043284659140218
$colJ[4,0] = 1.043308573552787
$colJ[5,0] = 1.043144112840303
$colJ[6,0] = 1.042456185870477
$colJ[7,0] = 1.041242127897351
$colJ[8,0] = 1.04043158180416
$colJ[9,0] = 1.040080341447529
$colJ[10,0] = 1.039949835861339
$colJ[11,0] = 1.039977831512886
$colJ[12,0] = 1.040069554610159
$colJ[13,0] = 1.040126063064768
$colJ[14,0] = 1.04045488690295
$colJ[15,0] = 1.040661078264535
$colJ[16,0] = 1.040781320214502
$colJ[17,0] = 1.040822315172506
$colJ[18,0] = 1.04063895856261
$colJ[19,0] = 1.040042545519859
$colJ[20,0] = 1.039667330242527
$colJ[21,0] = 1.039866260094643
$colJ[22,0] = 1.040648953584172
$colJ[23,0] = 1.041556202679001
$ws.Range("J2:J25").Value = $colJ

$colK = New-Object 'object[,]' 24,1
$colK[0,0] = 1.042130602937207
$colK[1,0] = 1.042588198571586
$colK[2,0] = 1.042883832314429
$colK[3,0] = 1.043008005224954
$colK[4,0] = 1.043028847831878
$colK[5,0] = 1.042885491958488
$colK[6,0] = 1.042285344366882
$colK[7,0] = 1.041224323500111
$colK[8,0] = 1.040514698543595
$colK[9,0] = 1.040206896407881
$colK[10,0] = 1.040092486456272
$colK[11,0] = 1.040117031315896
$colK[12,0] = 1.040197440841462
$colK[13,0] = 1.040246973429563
$colK[14,0] = 1.040535115285695
$colK[15,0] = 1.04071571795983
$colK[16,0] = 1.040821009249829
$colK[17,0] = 1.040856902142757
$colK[18,0] = 1.040696346289629
$colK[19,0] = 1.040173764399655
$colK[20,0] = 1.039844742583084
$colK[21,0] = 1.040019205935007
$colK[22,0] = 1.040705099671844
$colK[23,0] = 1.041499029033674
$ws.Range("K2:K25").Value = $colK

$colL = New-Object 'object[,]' 24,1
$colL[0,0] = 1.047616169536049
$colL[1,0] = 1.048243893789712
$colL[2,0] = 1.048650220296053
$colL[3,0] = 1.048821073916485
$colL[4,0] = 1.048849762946944
$colL[5,0] = 1.048652503118592
$colL[6,0] = 1.047828280727889
$colL[7,0] = 1.046377074858833
$colL[8,0] = 1.045410469421601
$colL[9,0] = 1.044992138893611
$colL[10,0] = 1.044836785959507
$colL[11,0] = 1.04487010814925
$colL[12,0] = 1.044979296678231
$colL[13,0] = 1.045046575829429
$colL[14,0] = 1.045438237294214
$colL[15,0] = 1.045683974859595
$colL[16,0] = 1.045827330184911
$colL[17,0] = 1.045876214131201
$colL[18,0] = 1.045657607388401
$colL[19,0] = 1.044947142439367
$colL[20,0] = 1.044500640210579
$colL[21,0] = 1.044737320579428
$colL[22,0] = 1.045669521649638
$colL[23,0] = 1.046752099165962
$ws.Range("L2:L25").Value = $colL

$colM = New-Object 'object[,]' 24,1
$colM[0,0] = 1.055745200838476
$colM[1,0] = 1.056535652856418
$colM[2,0] = 1.057047652390873
$colM[3,0] = 1.057263021126956
$colM[4,0] = 1.057299189724819
$colM[5,0] = 1.057050529672424
$colM[6,0] = 1.056012228432013
$colM[7,0] = 1.054186684373468
$colM[8,0] = 1.052972476047895
$colM[9,0] = 1.052447396459219
$colM[10,0] = 1.052252462119796
$colM[11,0] = 1.052294271487683
$colM[12,0] = 1.052431281004755
$colM[13,0] = 1.052515710858157
$colM[14,0] = 1.053007338305485
$colM[15,0] = 1.053315906385237
$colM[16,0] = 1.053495954506313
$colM[17,0] = 1.053557357345558
$colM[18,0] = 1.053282793157513
$colM[19,0] = 1.052390932246224
$colM[20,0] = 1.051830783950876
$colM[21,0] = 1.052127671833457
$colM[22,0] = 1.053297755397285
$colM[23,0] = 1.054658139976882
$ws.Range("M2:M25").Value = $colM

$ws.Range("N2").Value = 1.005712725503983

